# fix(gui) step 1 and 2
# Bump the quote date by one day and refresh the NEGRO / BLANCO price
# columns on the "SOPORTE DE ESTANTE" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quote date (A1) moves from 2024-01-17 to 2024-01-18
$ws.Range("A1").Value = 45309

# NEGRO price list (rows 26-32)
$ws.Range("D26").Value = 6797.942
$ws.Range("D27").Value = 8507.816999999999
$ws.Range("D28").Value = 11071.228
$ws.Range("D29").Value = 15643.851
$ws.Range("D30").Value = 23056.99
$ws.Range("D31").Value = 31703.345
$ws.Range("D32").Value = 40757.109

# BLANCO price list (rows 34-40) - mirrors the NEGRO prices above
$ws.Range("D34").Value = 6797.942
$ws.Range("D35").Value = 8507.816999999999
$ws.Range("D36").Value = 11071.228
$ws.Range("D37").Value = 15643.851
$ws.Range("D38").Value = 23056.99
$ws.Range("D39").Value = 31703.345
$ws.Range("D40").Value = 40757.109
